$wb = $excel.ActiveWorkbook

$listSheet   = $wb.Worksheets.Item("List")
$searchSheet = $wb.Worksheets.Item("Search")

# --- "List" sheet -----------------------------------------------------
# account.civility now renders its label instead of the raw enum
$listSheet.Range("F2").Value = '${account.civility.label}'

# account_addressId / account.addressId -> account_homeAddress / printer.print(account.homeAddress)
$listSheet.Range("K1").Value = '${msg.getProperty(' + "'" + 'account_homeAddress' + "'" + ')}'
$listSheet.Range("K2").Value = '${printer.print(account.homeAddress)}'

# --- "Search" sheet ----------------------------------------------------
# Insert a new criteria row for "Home Address" right before "Security Roles",
# pushing the existing Security Roles row (old row 14) down to row 15.
$searchSheet.Rows.Item(14).Insert()
$searchSheet.Range("A14").Value = '${msg.getProperty(' + "'" + 'account_homeAddress' + "'" + ')}'
$searchSheet.Range("B14").Value = '${homeAddress}'
